# Auto-generated edit script applying the Exodus_Profits workbook update
# Updates currentAveragePrice / LevePrice / LeveProfit columns across all 8 sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2015.3572
$ws.Range("J17").Value = 2015.3572
$ws.Range("L17").Value = 6046.071599999999
$ws.Range("N17").Value = -6382.071599999999
$ws.Range("H137").Value = 270698.97
$ws.Range("I137").Value = 2062.6
$ws.Range("J137").Value = 606494.4399999999
$ws.Range("K137").Value = 6187.799999999999
$ws.Range("L137").Value = 1819483.32
$ws.Range("M137").Value = -3637.799999999999
$ws.Range("N137").Value = -1824583.32

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6849.2144
$ws.Range("I32").Value = 3344.5122
$ws.Range("K32").Value = 3344.5122
$ws.Range("M32").Value = -3057.5122
$ws.Range("H64").Value = 48000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 48000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 48000
$ws.Range("N64").Value = -48496
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 48000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 48000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 48000
$ws.Range("N67").Value = -49716
$ws.Range("M67").ClearContents()
$ws.Range("H121").Value = 59924.547
$ws.Range("J121").Value = 59924.547
$ws.Range("L121").Value = 59924.547
$ws.Range("N121").Value = -63418.547
$ws.Range("H133").Value = 60000
$ws.Range("J133").Value = 70000
$ws.Range("L133").Value = 70000
$ws.Range("N133").Value = -75060
$ws.Range("H134").Value = 67017
$ws.Range("J134").Value = 67017
$ws.Range("L134").Value = 67017
$ws.Range("N134").Value = -77157
$ws.Range("H135").Value = 73600
$ws.Range("J135").Value = 73600
$ws.Range("L135").Value = 73600
$ws.Range("N135").Value = -83740
$ws.Range("H138").Value = 68563.39999999999
$ws.Range("J138").Value = 67106.75
$ws.Range("L138").Value = 67106.75
$ws.Range("N138").Value = -77386.75
$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360
$ws.Range("H141").Value = 92797
$ws.Range("J141").Value = 78661.664
$ws.Range("L141").Value = 78661.664
$ws.Range("N141").Value = -89021.664

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H56").Value = 10555
$ws.Range("J56").Value = 16110
$ws.Range("L56").Value = 16110
$ws.Range("N56").Value = -17588
$ws.Range("H80").Value = 544.5714
$ws.Range("I80").Value = 100
$ws.Range("K80").Value = 100
$ws.Range("M80").Value = 898
$ws.Range("H83").Value = 544.5714
$ws.Range("I83").Value = 100
$ws.Range("K83").Value = 500
$ws.Range("M83").Value = 4492

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 55000
$ws.Range("J63").Value = 55000
$ws.Range("L63").Value = 55000
$ws.Range("N63").Value = -56372
$ws.Range("H66").Value = 55000
$ws.Range("J66").Value = 55000
$ws.Range("L66").Value = 165000
$ws.Range("N66").Value = -171864
$ws.Range("H122").Value = 2313.8262
$ws.Range("I122").Value = 2071.3076
$ws.Range("K122").Value = 6213.9228
$ws.Range("M122").Value = -3763.9228

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 942.875
$ws.Range("I60").Value = 286.66666
$ws.Range("J60").Value = 1058.6765
$ws.Range("K60").Value = 859.9999799999999
$ws.Range("L60").Value = 3176.0295
$ws.Range("M60").Value = -608.9999799999999
$ws.Range("N60").Value = -3678.0295
$ws.Range("H88").Value = 4301
$ws.Range("I88").Value = 2000
$ws.Range("K88").Value = 6000
$ws.Range("M88").Value = -5572
$ws.Range("H91").Value = 4301
$ws.Range("I91").Value = 2000
$ws.Range("K91").Value = 6000
$ws.Range("M91").Value = -4518
$ws.Range("H138").Value = 6489.5576
$ws.Range("I138").Value = 6461.46
$ws.Range("K138").Value = 19384.38
$ws.Range("M138").Value = -14244.38

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 67500
$ws.Range("I64").Value = 67500
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 67500
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -67252
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 67500
$ws.Range("I67").Value = 67500
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 67500
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -66642
$ws.Range("N67").ClearContents()
$ws.Range("H122").Value = 3903.652
$ws.Range("I122").Value = 3549.2856
$ws.Range("J122").Value = 4454.8887
$ws.Range("K122").Value = 10647.8568
$ws.Range("L122").Value = 13364.6661
$ws.Range("M122").Value = -8197.856800000001
$ws.Range("N122").Value = -18264.6661
$ws.Range("H141").Value = 131993.5
$ws.Range("J141").Value = 149992
$ws.Range("L141").Value = 149992
$ws.Range("N141").Value = -160352

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3699.4
$ws.Range("J68").Value = 2999
$ws.Range("L68").Value = 2999
$ws.Range("N68").Value = -4497
$ws.Range("H71").Value = 3699.4
$ws.Range("J71").Value = 2999
$ws.Range("L71").Value = 14995
$ws.Range("N71").Value = -22483
$ws.Range("H101").Value = 46798.6
$ws.Range("J101").Value = 46798.6
$ws.Range("L101").Value = 46798.6
$ws.Range("N101").Value = -53288.6
$ws.Range("H122").Value = 66670828
$ws.Range("I122").Value = 66671190
$ws.Range("K122").Value = 200013570
$ws.Range("M122").Value = -200011120
$ws.Range("H134").Value = 131691.67
$ws.Range("J134").Value = 131691.67
$ws.Range("L134").Value = 131691.67
$ws.Range("N134").Value = -141831.67
$ws.Range("H135").Value = 76999
$ws.Range("J135").Value = 76999
$ws.Range("L135").Value = 76999
$ws.Range("N135").Value = -87139
$ws.Range("H138").Value = 108998.4
$ws.Range("J138").Value = 108998.4
$ws.Range("L138").Value = 108998.4
$ws.Range("N138").Value = -119278.4
$ws.Range("H140").Value = 69800
$ws.Range("J140").Value = 69800
$ws.Range("L140").Value = 69800
$ws.Range("N140").Value = -80160
$ws.Range("H141").Value = 78810
$ws.Range("J141").Value = 78810
$ws.Range("L141").Value = 78810
$ws.Range("N141").Value = -89170

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 98011
$ws.Range("J46").Value = 98011
$ws.Range("L46").Value = 98011
$ws.Range("N46").Value = -98473
$ws.Range("H86").Value = 54112.5
$ws.Range("J86").Value = 54112.5
$ws.Range("L86").Value = 54112.5
$ws.Range("N86").Value = -56358.5
$ws.Range("H89").Value = 54112.5
$ws.Range("J89").Value = 54112.5
$ws.Range("L89").Value = 270562.5
$ws.Range("N89").Value = -281794.5
$ws.Range("H133").Value = 66748.25
$ws.Range("J133").Value = 63831
$ws.Range("L133").Value = 63831
$ws.Range("N133").Value = -73951
$ws.Range("H134").Value = 98011
$ws.Range("J134").Value = 98011
$ws.Range("L134").Value = 294033
$ws.Range("N134").Value = -299103
$ws.Range("H137").Value = 138249.5
$ws.Range("J137").Value = 138249.5
$ws.Range("L137").Value = 138249.5
$ws.Range("N137").Value = -148449.5
$ws.Range("H140").Value = 150000
$ws.Range("J140").Value = 150000
$ws.Range("L140").Value = 150000
$ws.Range("N140").Value = -160360
$ws.Range("H141").Value = 61064.375
$ws.Range("J141").Value = 61064.375
$ws.Range("L141").Value = 61064.375
$ws.Range("N141").Value = -71424.375
